$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing first table (rows 1-8) with new values ---
$ws.Range("B2").Value = 2850
$ws.Range("C2").Value = 410

$ws.Range("B3").Value = 6300
$ws.Range("C3").Value = 365

$ws.Range("B4").Value = 95000
$ws.Range("C4").Value = 300

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 10

# --- Add the new second table (rows 11-18): Cryogenic engines support ---

# Header row 11
$ws.Range("A11").Value = "Methalox Variant"
$ws.Range("C11").Value = "ISP"
$ws.Range("G11").Value = "Hydrolox ISP"
$ws.Range("H11").Value = "Extended"

# Row 12
$ws.Range("A12").Value = "Engine"
$ws.Range("B12").Value = 360
$ws.Range("C12").Formula = "=0.82*G12"
$ws.Range("D12").Formula = "=B12*1.2"
$ws.Range("E12").Formula = "=C12*1.1"
$ws.Range("G12").Formula = "=H12/1.1"
$ws.Range("H12").Value = 514.8

# Row 13
$ws.Range("A13").Value = "ASL Thrust"
$ws.Range("B13").Formula = "=(C13/C12)*B12"
$ws.Range("D13").Formula = "=(E13/E12)*D12"
$ws.Range("H13").Value = 137.5

# Rows 14-15 weight values
$ws.Range("H14").Value = 27.5
$ws.Range("H15").Value = 11

# Shared formula groups (fill down like Excel does to produce shared formulas)
$ws.Range("C13:C18").Formula = "=0.82*G13"
$ws.Range("E13:E18").Formula = "=C13*1.1"
$ws.Range("G13:G15").Formula = "=H13/1.1"

# --- Selection as in target file ---
$ws.Range("B14").Select()
